# Daily attendance processing - 2026-01-31 20:40:01
#
# Normalizes the "Recorded By" (column G) entries: whenever the comma
# separated list of recorders has its two "interesting" entries out of
# order (the real submitter email vs. the "System"/duplicate-system
# marker), the last two entries are swapped into the corrected order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Test-NeedsSwap($s) {
    if ($s -eq $null) { return $false }
    $parts = $s -split ", "
    if ($parts.Length -lt 2) { return $false }

    $hasDnasr = $false
    foreach ($p in $parts) {
        if ($p -eq "dnasr281@gmail.com") { $hasDnasr = $true }
    }

    $hasUpperSystem = $false
    $hasLowerSystem = $false
    foreach ($p in $parts) {
        if ($p.Equals("System")) { $hasUpperSystem = $true }
        if ($p.Equals("system")) { $hasLowerSystem = $true }
    }
    $hasDupSystem = $hasUpperSystem -and $hasLowerSystem

    return ($hasDnasr -or $hasDupSystem)
}

function Get-SwappedLastTwo($s) {
    $parts = $s -split ", "
    $n = $parts.Length
    $tmp = $parts[$n - 1]
    $parts[$n - 1] = $parts[$n - 2]
    $parts[$n - 2] = $tmp
    return ($parts -join ", ")
}

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value()
    if (Test-NeedsSwap $val) {
        $cell.Value = Get-SwappedLastTwo $val
    }
}
